$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    3 = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    4 = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.611132179096228)
    5 = @(0.01293466051926884, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.30063094841924)
    6 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    7 = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1.672833113781282)
    8 = @(0.04271373187048222, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.330410019770453)
    9 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
